# Apply updated coin-price snapshot values scraped on 2023-02-05.
# Numeric-looking text (prices / percentages) is stored as text in the
# sheet, so values are written with a leading "'" quote-prefix to stop
# Excel from auto-converting them to numbers (preserves original text +
# General number format / style, matching the source data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'335.60"
$ws.Range("E2").Value = "'1.79%"

# Row 3
$ws.Range("D3").Value = "'44.02"
$ws.Range("E3").Value = "'6.23%"

# Row 4
$ws.Range("D4").Value = "'5.761"
$ws.Range("E4").Value = "'2.05%"

# Row 5
$ws.Range("D5").Value = "'0.08393"
$ws.Range("E5").Value = "'1.36%"

# Row 6
$ws.Range("D6").Value = "'8.861"
$ws.Range("E6").Value = "'1.12%"

# Row 7
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.959"
$ws.Range("E7").Value = "'-4.05%"

# Row 8
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.848"
$ws.Range("E8").Value = "'-4.13%"

# Row 9
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9477"
$ws.Range("E9").Value = "'2.43%"

# Row 10
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1239"
$ws.Range("E10").Value = "'-2.82%"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1982"
$ws.Range("E11").Value = "'1.09%"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.1014"
$ws.Range("E12").Value = "'7.59%"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04423"
$ws.Range("E13").Value = "'11.34%"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1068"
$ws.Range("E14").Value = "'0.69%"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001291"
$ws.Range("E15").Value = "'-1.38%"

# Row 16
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006085"
$ws.Range("E16").Value = "'-0.24%"

# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.495"
$ws.Range("E17").Value = "'1.47%"

# Row 18
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.530"
$ws.Range("E18").Value = "'-0.06%"

# Row 20
$ws.Range("D20").Value = "'8.656"
$ws.Range("E20").Value = "'3.41%"

# Row 21
$ws.Range("E21").Value = "'-0.68%"

# Row 22
$ws.Range("D22").Value = "'0.2642"
$ws.Range("E22").Value = "'-0.76%"

# Row 23
$ws.Range("D23").Value = "'0.04419"
$ws.Range("E23").Value = "'0.54%"

# Row 24
$ws.Range("D24").Value = "'0.001257"
$ws.Range("E24").Value = "'0.11%"

# Row 25
$ws.Range("D25").Value = "'0.004357"
$ws.Range("E25").Value = "'0.94%"

# Row 26
$ws.Range("E26").Value = "'5.12%"

# Row 27
$ws.Range("D27").Value = "'0.0003999"
$ws.Range("E27").Value = "'-94.67%"

# Row 39
$ws.Range("D39").Value = "'0.02840"
$ws.Range("E39").Value = "'2.57%"

# Row 40
$ws.Range("D40").Value = "'0.05895"
$ws.Range("E40").Value = "'7.00%"

# Row 41
$ws.Range("D41").Value = "'0.007914"
$ws.Range("E41").Value = "'-0.03%"

# Row 42
$ws.Range("D42").Value = "'0.1429"
$ws.Range("E42").Value = "'0.55%"

# Row 43
$ws.Range("D43").Value = "'0.009014"
$ws.Range("E43").Value = "'0.80%"

# Row 44
$ws.Range("D44").Value = "'0.002146"
$ws.Range("E44").Value = "'0.22%"

# Row 45
$ws.Range("D45").Value = "'0.009908"
$ws.Range("E45").Value = "'-16.14%"

# Row 46
$ws.Range("D46").Value = "'0.00007235"
$ws.Range("E46").Value = "'3.18%"

# Row 47
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.19%"

# Row 48
$ws.Range("D48").Value = "'0.003203"

# Row 49
$ws.Range("D49").Value = "'0.002276"
$ws.Range("E49").Value = "'-0.20%"

# Row 50
$ws.Range("D50").Value = "'0.00002106"
$ws.Range("E50").Value = "'0.19%"

# Row 51
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'0.19%"
